$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.803.29'
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').Value = '3.917.15'
$ws.Range('E3').Value = '  +2.93%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '603.74'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '167.60'
$ws.Range('E6').Value = '  +2.26%  '
$ws.Range('D7').Value = '3.912.11'
$ws.Range('E7').Value = '  +2.88%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.532'
$ws.Range('E9').Value = '  -0.39%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.169'
$ws.Range('E10').Value = '  +0.14%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.48'
$ws.Range('E11').Value = '  +2.76%  '
$ws.Range('E12').Value = '  +1.02%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000254'
$ws.Range('E13').Value = '  +3.59%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '37.54'
$ws.Range('E14').Value = '  +1.23%  '
$ws.Range('D15').Value = '4.576.74'
$ws.Range('E15').Value = '  +3.08%  '
$ws.Range('D16').Value = '3.899.54'
$ws.Range('E16').Value = '  +2.18%  '
$ws.Range('D17').Value = '68.900.32'
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '7.49'
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '17.39'
$ws.Range('E19').Value = '  +0.99%  '
$ws.Range('E20').Value = '  -1.96%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.06'
$ws.Range('E21').Value = '  -3.25%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '492.96'
$ws.Range('E22').Value = '  +1.20%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.729'
$ws.Range('E23').Value = '  +1.43%  '
$ws.Range('E24').Value = '  +4.09%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '84.84'
$ws.Range('E25').Value = '  +0.27%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E27').Value = '  -0.34%  '
$ws.Range('E29').Value = '  +0.13%  '
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('D31').Value = '4.071.54'
$ws.Range('E31').Value = '  +2.76%  '
$ws.Range('E32').Value = '  -0.26%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '7.75'
$ws.Range('E33').Value = '  -3.06%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '31.93'
$ws.Range('E34').Value = '  +0.43%  '
$ws.Range('D35').Value = '3.879.89'
$ws.Range('E35').Value = '  +3.47%  '
$ws.Range('E36').Value = '  +0.23%  '
$ws.Range('E37').Value = '  +1.34%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.140'
$ws.Range('E38').Value = '  +0.32%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.96'
$ws.Range('E39').Value = '  +1.65%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.26'
$ws.Range('E40').Value = '  +7.22%  '
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('E42').Value = '  +0.08%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '437.33'
$ws.Range('E43').Value = '  -0.15%  '
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '47.96'
$ws.Range('E45').Value = '  -1.26%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '8.58'
$ws.Range('E46').Value = '  +2.69%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '143.04'
$ws.Range('E48').Value = '  +0.79%  '
$ws.Range('D49').Value = '2.820.56'
$ws.Range('E49').Value = '  -0.21%  '
$ws.Range('E50').Value = '  +18.10%  '
$ws.Range('E51').Value = '  +1.81%  '
